$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add names to column B, rows 2-5 (these become shared strings: mahi, dom, eric, john)
$ws.Range("B2").Value = "mahi"
$ws.Range("B3").Value = "dom"
$ws.Range("B4").Value = "eric"
$ws.Range("B5").Value = "john"

# Update the active cell/selection to D11, matching the diff
$ws.Range("D11").Select()
